$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 244 entirely; everything below shifts up by one row.
$ws.Rows("244").Delete()
